# Update week 8 documents.
#
# Two textual edits are made:
#   1. "FirstNameLastNameWeek7" -> "FirstNameLastNameWeek8" (the digit is
#      changed from 7 to 8, which - because of how the edit originally
#      happened in Word - leaves the surrounding sentence split into three
#      runs instead of one, with the replaced "8" living in its own run).
#   2. The word "two-tier" is split into "t" / "wo-tier" at the point where
#      the cursor was left last, which is exactly where Word parks its
#      "_GoBack" bookmark when the document is saved. Re-adding a bookmark
#      named "_GoBack" automatically relocates it (bookmark names are
#      unique), so the old "_GoBack" bookmark that used to sit after
#      "Hiding data" disappears as a side effect, matching the target.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Week7" -> "Week8", with "8" ending up as its own run.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Week7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    # $rng1 now spans the matched text "Week7"; the final character is the
    # "7" that needs to become "8".
    $digitStart = $rng1.End - 1
    $digitEnd = $rng1.End

    # Temporarily bookmark the two collapsed points immediately around the
    # digit. This forces Word to split the run there without disturbing any
    # of the other runs in the paragraph.
    $splitBefore = $d.Range($digitStart, $digitStart)
    $d.Bookmarks.Add("TempSplitBefore", $splitBefore)
    $splitAfter = $d.Range($digitEnd, $digitEnd)
    $d.Bookmarks.Add("TempSplitAfter", $splitAfter)

    # The digit is now isolated in its own run; replace it.
    $digitRange = $d.Range($digitStart, $digitEnd)
    $digitRange.Text = "8"

    # Drop the helper bookmarks - the run split they created remains intact.
    $d.Bookmarks("TempSplitBefore").Delete()
    $d.Bookmarks("TempSplitAfter").Delete()
}

# ---------------------------------------------------------------------
# Edit 2: split "two-tier" into "t" | "wo-tier" and drop a "_GoBack"
# bookmark at the split point (this also removes the old "_GoBack"
# bookmark wherever it used to be, since bookmark names are unique).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("two-tier", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    # Split point is right after the first character ("t") of "two-tier".
    $splitPoint = $d.Range($rng2.Start + 1, $rng2.Start + 1)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}
